# Auto update stock data
# Updates the "as of" date (column A) for each ticker's most-recent-row
# from 2025/11/08 to 2025/11/09, keeping the value as literal text
# (not converting it into a real Excel date serial / changing cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    # Leading apostrophe forces the value to be stored as text instead of
    # being auto-parsed into a date value by the input parser.
    $cell.Value = "'2025/11/09"
    # Drop the "quote prefix" formatting that the apostrophe trick applies,
    # restoring the cell to its original (unstyled) appearance while
    # keeping the text content intact.
    $cell.ClearFormats()
}
